# New crime data collected — update weekly CompStat report (70th Precinct)
# for the week of 12/5/2022 through 12/11/2022 (Volume 29, Number 49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (report title / reporting-week banner)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# ---------------------------------------------------------------------
# Helper: make a cell hold the literal text "0" / "***.*" the same way
# the rest of the sheet represents "insufficient data" cells (style 14,
# shared-string text), instead of letting it auto-convert to a number.
# We enter it as text (leading apostrophe) then copy the number format
# from an existing "style 14" text cell (column C of row 23, the
# template "all N/A" row) onto it so the stored style index matches.
# ---------------------------------------------------------------------
$naTemplate = "C23"

function Set-TextCell($ws, $addr, $text, $naTemplate) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($naTemplate).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-NumericCell($ws, $addr, $value, $templateAddr) {
    $ws.Range($addr).Value = $value
    if ($templateAddr) {
        $ws.Range($templateAddr).Copy() | Out-Null
        $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    }
}

# ---------------------------------------------------------------------
# Row 15 — Rape
# ---------------------------------------------------------------------
$ws.Range("C15").Value = 1
Set-TextCell $ws "D15" "0" $naTemplate
Set-TextCell $ws "E15" "***.*" $naTemplate
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 30
$ws.Range("J15").Value = 24
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = -72.477064220183

# ---------------------------------------------------------------------
# Row 16 — Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 197
$ws.Range("J16").Value = 168
$ws.Range("K16").Value = 17.261904761904
$ws.Range("L16").Value = 15.204678362573
$ws.Range("M16").Value = -48.697916666666
$ws.Range("N16").Value = -90.473887814313

# ---------------------------------------------------------------------
# Row 17 — Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -72.727272727272
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 338
$ws.Range("J17").Value = 336
$ws.Range("K17").Value = 0.595238095238
$ws.Range("L17").Value = 7.643312101910
$ws.Range("M17").Value = -3.703703703703
$ws.Range("N17").Value = -59.129383313180

# ---------------------------------------------------------------------
# Row 18 — Burglary
# ---------------------------------------------------------------------
Set-NumericCell $ws "C18" 6 "D18"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -47.619047619047
$ws.Range("I18").Value = 180
$ws.Range("J18").Value = 197
$ws.Range("K18").Value = -8.629441624365
$ws.Range("L18").Value = -12.195121951219
$ws.Range("M18").Value = -41.558441558441
$ws.Range("N18").Value = -93.433053630062

# ---------------------------------------------------------------------
# Row 19 — Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -33.898305084745
$ws.Range("I19").Value = 489
$ws.Range("J19").Value = 480
$ws.Range("K19").Value = 1.875
$ws.Range("L19").Value = 7.947019867549
$ws.Range("M19").Value = -22.626582278481
$ws.Range("N19").Value = -54.213483146067

# ---------------------------------------------------------------------
# Row 20 — G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 55.555555555555
$ws.Range("I20").Value = 115
$ws.Range("J20").Value = 121
$ws.Range("K20").Value = -4.958677685950
$ws.Range("L20").Value = -24.342105263157
$ws.Range("M20").Value = -36.111111111111
$ws.Range("N20").Value = -94.967177242888

# ---------------------------------------------------------------------
# Row 21 — TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -21.212121212121
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 140
$ws.Range("H21").Value = -22.857142857142
$ws.Range("I21").Value = 1351
$ws.Range("J21").Value = 1329
$ws.Range("K21").Value = 1.655379984951
$ws.Range("L21").Value = 1.274362818590
$ws.Range("M21").Value = -28.366914103923
$ws.Range("N21").Value = -85.222052067381

# ---------------------------------------------------------------------
# Row 22 — Transit
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 1
Set-TextCell $ws "D22" "0" $naTemplate
Set-TextCell $ws "E22" "***.*" $naTemplate
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = -27.777777777777
$ws.Range("L22").Value = -13.333333333333
$ws.Range("M22").Value = -48

# Row 23 — Housing: unchanged (all "N/A" placeholders already in place).

# ---------------------------------------------------------------------
# Row 24 — Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 79
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 97.5
$ws.Range("F24").Value = 179
$ws.Range("G24").Value = 123
$ws.Range("H24").Value = 45.528455284552
$ws.Range("I24").Value = 1472
$ws.Range("J24").Value = 1317
$ws.Range("K24").Value = 11.769172361427
$ws.Range("L24").Value = 10.926902788244
$ws.Range("M24").Value = 14.463452566096

# ---------------------------------------------------------------------
# Row 25 — Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -37.5
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -14
$ws.Range("I25").Value = 537
$ws.Range("J25").Value = 541
$ws.Range("K25").Value = -0.739371534195
$ws.Range("L25").Value = 0.750469043151
$ws.Range("M25").Value = -27.822580645161

# ---------------------------------------------------------------------
# Row 26 — UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 2
Set-TextCell $ws "D26" "0" $naTemplate
Set-TextCell $ws "E26" "***.*" $naTemplate
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 46
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 15
$ws.Range("L26").Value = 9.523809523809

# ---------------------------------------------------------------------
# Row 27 — Other Sex Crimes
# ---------------------------------------------------------------------
Set-NumericCell $ws "C27" 1 "C17"
Set-NumericCell $ws "D27" 5 "D17"
Set-NumericCell $ws "E27" -80 "E17"
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("J27").Value = 62
$ws.Range("K27").Value = 16.129032258064
$ws.Range("L27").Value = 4.347826086956

# ---------------------------------------------------------------------
# Row 28 — Shooting Vic. (only the 2-Year % column changed)
# ---------------------------------------------------------------------
$ws.Range("L28").Value = -68.181818181818

# ---------------------------------------------------------------------
# Row 29 — Shooting Inc. (only the 2-Year % column changed)
# ---------------------------------------------------------------------
$ws.Range("L29").Value = -58.064516129032
